{"js": "// Replacements, listed in document order. A couple of source strings\n// (e.g. \"92\u00f79=\") repeat more than once in the document, so we search\n// once per distinct text and then consume the N-th match (in document\n// order) for the N-th time that text is requested.\nconst replacements = [\n  [\"2025-08-31 Sunday\", \"2025-09-01 Monday\"],\n  [\"33\u00f72=\", \"19\u00f79=\"],\n  [\"42\u00f77=\", \"20\u00f75=\"],\n  [\"10\u00f75=\", \"62\u00f79=\"],\n  [\"89\u00f74=\", \"97\u00f77=\"],\n  [\"62\u00f77=\", \"24\u00f78=\"],\n  [\"66\u00f76=\", \"70\u00f72=\"],\n  [\"92\u00f79=\", \"82\u00f72=\"],\n  [\"57\u00f75=\", \"16\u00f78=\"],\n  [\"50\u00f73=\", \"61\u00f79=\"],\n  [\"15\u00f78=\", \"85\u00f77=\"],\n  [\"80\u00f75=\", \"79\u00f75=\"],\n  [\"96\u00f77=\", \"16\u00f75=\"],\n  [\"82\u00f75=\", \"30\u00f73=\"],\n  [\"87\u00f77=\", \"78\u00f78=\"],\n  [\"48\u00f72=\", \"97\u00f75=\"],\n  [\"11\u00f72=\", \"80\u00f79=\"],\n  [\"59\u00f72=\", \"93\u00f76=\"],\n  [\"57\u00f77=\", \"51\u00f73=\"],\n  [\"92\u00f79=\", \"91\u00f76=\"],\n  [\"98\u00f75=\", \"42\u00f74=\"],\n  [\"50\u00f76=\", \"98\u00f74=\"],\n  [\"60\u00f72=\", \"34\u00f76=\"],\n  [\"74\u00f74=\", \"73\u00f78=\"],\n  [\"28\u00f74=\", \"63\u00f73=\"],\n  [\"56\u00f72=\", \"81\u00f78=\"],\n];\n\n// Search once per distinct source text, loading the matching ranges.\nconst searchResults = {};\nfor (const [findText] of replacements) {\n  if (!(findText in searchResults)) {\n    const results = context.document.body.search(findText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    searchResults[findText] = results;\n  }\n}\nawait context.sync();\n\n// Walk the replacement list again (document order), and for each entry\n// take the next not-yet-used match for that text.\nconst nextIndex = {};\nfor (const [findText, replaceText] of replacements) {\n  const results = searchResults[findText];\n  const i = nextIndex[findText] || 0;\n  results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  nextIndex[findText] = i + 1;\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction ReplaceOnce($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1) | Out-Null\n}\n\n# Replacements are applied in document order so duplicate source values\n# (e.g. \"92\u00f79=\" appearing twice) map to the correct distinct targets.\nReplaceOnce \"2025-08-31 Sunday\" \"2025-09-01 Monday\"\nReplaceOnce \"33\u00f72=\" \"19\u00f79=\"\nReplaceOnce \"42\u00f77=\" \"20\u00f75=\"\nReplaceOnce \"10\u00f75=\" \"62\u00f79=\"\nReplaceOnce \"89\u00f74=\" \"97\u00f77=\"\nReplaceOnce \"62\u00f77=\" \"24\u00f78=\"\nReplaceOnce \"66\u00f76=\" \"70\u00f72=\"\nReplaceOnce \"92\u00f79=\" \"82\u00f72=\"\nReplaceOnce \"57\u00f75=\" \"16\u00f78=\"\nReplaceOnce \"50\u00f73=\" \"61\u00f79=\"\nReplaceOnce \"15\u00f78=\" \"85\u00f77=\"\nReplaceOnce \"80\u00f75=\" \"79\u00f75=\"\nReplaceOnce \"96\u00f77=\" \"16\u00f75=\"\nReplaceOnce \"82\u00f75=\" \"30\u00f73=\"\nReplaceOnce \"87\u00f77=\" \"78\u00f78=\"\nReplaceOnce \"48\u00f72=\" \"97\u00f75=\"\nReplaceOnce \"11\u00f72=\" \"80\u00f79=\"\nReplaceOnce \"59\u00f72=\" \"93\u00f76=\"\nReplaceOnce \"57\u00f77=\" \"51\u00f73=\"\nReplaceOnce \"92\u00f79=\" \"91\u00f76=\"\nReplaceOnce \"98\u00f75=\" \"42\u00f74=\"\nReplaceOnce \"50\u00f76=\" \"98\u00f74=\"\nReplaceOnce \"60\u00f72=\" \"34\u00f76=\"\nReplaceOnce \"74\u00f74=\" \"73\u00f78=\"\nReplaceOnce \"28\u00f74=\" \"63\u00f73=\"\nReplaceOnce \"56\u00f72=\" \"81\u00f78=\"\n"}
